$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'92.543.85"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.30%  '

$ws.Range('D3').Value = "'3.095.90"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.70%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = "'239.58"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '

$ws.Range('D6').Value = "'608.87"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.41%  '

$ws.Range('D7').Value = "'1.10"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.86%  '

$ws.Range('D8').Value = "'0.388"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.16%  '

$ws.Range('D9').Value = "'1.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.20%  '

$ws.Range('D10').Value = "'3.091.41"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.78%  '

$ws.Range('D11').Value = "'0.750"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.61%  '

$ws.Range('D12').Value = "'0.200"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.28%  '

$ws.Range('D13').Value = "'0.0000247"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.92%  '

$ws.Range('D14').Value = "'92.575.83"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('D15').Value = "'34.03"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.12%  '

$ws.Range('D16').Value = "'5.42"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.14%  '

$ws.Range('D17').Value = "'3.686.86"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.17%  '

$ws.Range('D18').Value = "'3.108.44"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').Value = "'3.75"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.40%  '

$ws.Range('D20').Value = "'14.67"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.98%  '

$ws.Range('D21').Value = "'5.73"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.97%  '

$ws.Range('D22').Value = "'443.71"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.50%  '

$ws.Range('D23').Value = "'9.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.22%  '

$ws.Range('D24').Value = "'0.0000199"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.36%  '

$ws.Range('D25').Value = "'5.71"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.18%  '

$ws.Range('D26').Value = "'86.26"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.05%  '

$ws.Range('D27').Value = "'11.58"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.63%  '

$ws.Range('D28').Value = "'3.274.57"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.30%  '

$ws.Range('D29').Value = "'0.997"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.26%  '

$ws.Range('E30').Value = '  +11.10%  '

$ws.Range('D31').Value = "'0.231"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('D32').Value = "'0.168"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.61%  '

$ws.Range('D33').Value = "'9.08"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.38%  '

$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = "'0.995"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.64%  '

$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D35').Value = "'7.98"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.97%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.158"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.37%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = "'25.85"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.60%  '

$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').Value = "'1.89"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.14%  '

$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = "'486.08"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.95%  '

$ws.Range('E40').Value = '  +1.74%  '

$ws.Range('D41').Value = "'1.28"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.42%  '

$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = "'23.53"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.15%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = "'0.428"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.03%  '

$ws.Range('D44').Value = "'3.33"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.76%  '

$ws.Range('E45').Value = '  +0.03%  '

$ws.Range('D46').Value = "'162.79"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.79%  '

$ws.Range('D47').Value = "'1.87"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.97%  '

$ws.Range('D48').Value = "'0.680"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.84%  '

$ws.Range('D49').Value = "'1.37"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.07%  '

$ws.Range('D50').Value = "'0.0329"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.92%  '

$ws.Range('D51').Value = "'43.99"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.16%  '
